$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Heure de départ" (H) and "Tâches réalisées" (I) for the second
# table (F2:I6) - these were previously left blank. Text is entered in the
# same order it was typed by the author (row 3, 5, 6, then 4) so that the
# shared-string table is populated in the matching sequence.
$ws.Range("H3").Value = 0.72222222222222221
$ws.Range("H3").NumberFormat = "h:mm"
$ws.Range("I3").Value = "Schéma électrique, début du programme Arduino"

$ws.Range("H5").Value = 0.72916666666666663
$ws.Range("H5").NumberFormat = "h:mm"
$ws.Range("I5").Value = "Programme en C"

$ws.Range("H6").Value = 0.75
$ws.Range("H6").NumberFormat = "h:mm"
$ws.Range("I6").Value = "Gantt, Trello, Rapport, Tests sur Processing + Arduino"

$ws.Range("H4").Value = 0.75
$ws.Range("H4").NumberFormat = "h:mm"
$ws.Range("I4").Value = "Continue le programme C (Structures, commentaires, nombre de lignes en mémoire)"

# Widen column I to fit the new, much longer "Tâches réalisées" text.
$ws.Columns.Item(9).ColumnWidth = 67.5

# Update the active selection to reflect where the user ended up after
# entering this data.
$ws.Range("H7").Select()
